# Insert a new weekly data row for "Terminal Hortofrutícola Agro Chillán - Piña"
# (commit: "Fruta / hortaliza, semanal") at row 298, pushing the existing
# rows 298-377 down to 299-378.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 298 (shifts 298:377 -> 299:378)
$ws.Rows("298").Insert()

# Populate the newly inserted row with the new observation
$ws.Cells.Item(298, 1).Value = 7
$ws.Cells.Item(298, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(298, 3).Value = "Ñuble"
$ws.Cells.Item(298, 4).Value = 45204
$ws.Cells.Item(298, 5).Value = 16
$ws.Cells.Item(298, 6).Value = "Fruta"
$ws.Cells.Item(298, 7).Value = 100108
$ws.Cells.Item(298, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(298, 9).Value = 100108005
$ws.Cells.Item(298, 10).Value = "Piña"
$ws.Cells.Item(298, 11).Value = "Caramelo"
$ws.Cells.Item(298, 12).Value = "Segunda"
$ws.Cells.Item(298, 13).Value = 120
$ws.Cells.Item(298, 14).Value = 21000
$ws.Cells.Item(298, 15).Value = 22000
$ws.Cells.Item(298, 16).Value = 21500
$ws.Cells.Item(298, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(298, 18).Value = "Ecuador"
$ws.Cells.Item(298, 19).Value = 1536
$ws.Cells.Item(298, 20).Value = 14
